$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column B (rows 2-23)
$bValues = @(
    "NSE:CARYSIL",
    "NSE:CYIENTDLM",
    "NSE:EGOLD",
    "NSE:FEDERALBNK",
    "NSE:GHCL",
    "NSE:GMBREW",
    "NSE:GRINFRA",
    "NSE:HDFCGOLD",
    "NSE:HINDALCO",
    "NSE:IIFL",
    "NSE:IMFA",
    "NSE:KEYFINSERV",
    "NSE:KSL",
    "NSE:LATENTVIEW",
    "NSE:LATTEYS",
    "NSE:MMFL",
    "NSE:MSPL",
    "NSE:NAGREEKCAP",
    "NSE:NEOGEN",
    "NSE:NESTLEIND",
    "NSE:ONEPOINT",
    "NSE:RAMCOSYS"
)

# New values for column C (rows 2-14); rows 15-23 become empty
$cValues = @(
    "NSE:ARVIND",
    "NSE:ASTRON",
    "NSE:CENTRALBK",
    "NSE:CUB",
    "NSE:DHANUKA",
    "NSE:DIAMINESQ",
    "NSE:DYNPRO",
    "NSE:FACT",
    "NSE:HITECHCORP",
    "NSE:KOTHARIPRO",
    "NSE:LFIC",
    "NSE:MAHABANK",
    "NSE:PRICOLLTD"
)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

for ($i = 0; $i -lt $cValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
}

# Clear column C for rows 15-23 (formerly had values, now empty)
for ($row = 15; $row -le 23; $row++) {
    $ws.Cells.Item($row, 3).Value = ""
}

# Delete rows 24-32 entirely (reduces dimension to A1:F23)
$ws.Range("A24:F32").EntireRow.Delete()
